$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "64.068.93"
Set-TextValue $ws.Range("E2") "  +0.69%  "
Set-TextValue $ws.Range("D3") "2.757.24"
Set-TextValue $ws.Range("E3") "  +1.19%  "
Set-TextValue $ws.Range("E4") "  -0.11%  "
Set-TextValue $ws.Range("D5") "578.23"
Set-TextValue $ws.Range("E5") "  +0.37%  "
Set-TextValue $ws.Range("D6") "158.39"
Set-TextValue $ws.Range("E6") "  +2.65%  "
Set-TextValue $ws.Range("E7") "  +0.25%  "
Set-TextValue $ws.Range("E8") "  +0.08%  "
Set-TextValue $ws.Range("E9") "  -1.63%  "
Set-TextValue $ws.Range("E10") "  -14.70%  "
Set-TextValue $ws.Range("D11") "0.386"
Set-TextValue $ws.Range("E11") "  -1.01%  "
Set-TextValue $ws.Range("E12") "  -2.26%  "
Set-TextValue $ws.Range("D13") "3.245.85"
Set-TextValue $ws.Range("E13") "  +1.11%  "
Set-TextValue $ws.Range("E14") "  +2.41%  "
Set-TextValue $ws.Range("D15") "63.749.64"
Set-TextValue $ws.Range("E15") "  +0.33%  "
Set-TextValue $ws.Range("E16") "  -0.15%  "
Set-TextValue $ws.Range("D17") "2.761.09"
Set-TextValue $ws.Range("D18") "12.18"
Set-TextValue $ws.Range("E18") "  +1.97%  "
Set-TextValue $ws.Range("D19") "4.88"
Set-TextValue $ws.Range("E19") "  +0.42%  "
Set-TextValue $ws.Range("D20") "360.00"
Set-TextValue $ws.Range("E20") "  -0.09%  "
Set-TextValue $ws.Range("D21") "6.84"
Set-TextValue $ws.Range("E21") "  -1.29%  "
Set-TextValue $ws.Range("D22") "0.548"
Set-TextValue $ws.Range("E22") "  +2.67%  "
Set-TextValue $ws.Range("D23") "0.999"
Set-TextValue $ws.Range("E23") "  +0.27%  "
Set-TextValue $ws.Range("D24") "65.77"
Set-TextValue $ws.Range("E24") "  -0.13%  "
Set-TextValue $ws.Range("E25") "  +1.45%  "
Set-TextValue $ws.Range("D26") "8.53"
Set-TextValue $ws.Range("E26") "  +0.07%  "
Set-TextValue $ws.Range("D27") "0.997"
Set-TextValue $ws.Range("E27") "  +0.04%  "
Set-TextValue $ws.Range("D28") "0.0₃0925"
Set-TextValue $ws.Range("E28") "  +2.10%  "
Set-TextValue $ws.Range("D29") "1.96"
Set-TextValue $ws.Range("E29") "  -1.31%  "
Set-TextValue $ws.Range("E30") "  -1.13%  "
Set-TextValue $ws.Range("E31") "  +0.99%  "
Set-TextValue $ws.Range("D32") "167.46"
Set-TextValue $ws.Range("E32") "  -2.46%  "
Set-TextValue $ws.Range("D33") "20.33"
Set-TextValue $ws.Range("E33") "  -0.67%  "
Set-TextValue $ws.Range("D34") "4.95"
Set-TextValue $ws.Range("E34") "  +3.67%  "
Set-TextValue $ws.Range("E35") "  +0.14%  "
Set-TextValue $ws.Range("E36") "  +1.91%  "
Set-TextValue $ws.Range("E37") "  +0.12%  "
Set-TextValue $ws.Range("D38") "0.991"
Set-TextValue $ws.Range("E38") "  -0.67%  "
Set-TextValue $ws.Range("D39") "6.26"
Set-TextValue $ws.Range("E39") "  +11.92%  "
Set-TextValue $ws.Range("E40") "  -1.19%  "
Set-TextValue $ws.Range("D41") "329.65"
Set-TextValue $ws.Range("E41") "  -4.33%  "
Set-TextValue $ws.Range("D42") "39.39"
Set-TextValue $ws.Range("E42") "  +0.39%  "
Set-TextValue $ws.Range("D43") "21.59"
Set-TextValue $ws.Range("E43") "  -0.42%  "
Set-TextValue $ws.Range("D44") "0.0595"
Set-TextValue $ws.Range("E44") "  +0.76%  "
Set-TextValue $ws.Range("E45") "  +0.07%  "
Set-TextValue $ws.Range("D46") "0.0257"
Set-TextValue $ws.Range("E46") "  +1.20%  "
Set-TextValue $ws.Range("D47") "0.636"
Set-TextValue $ws.Range("E47") "  -1.09%  "
Set-TextValue $ws.Range("D48") "136.64"
Set-TextValue $ws.Range("E48") "  -1.96%  "
Set-TextValue $ws.Range("E49") "  +0.83%  "
Set-TextValue $ws.Range("E50") "  +0.45%  "
Set-TextValue $ws.Range("E51") "  +0.73%  "
